$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 822, shifting rows 822:863 down to 823:864
$ws.Rows.Item(822).Insert()

$ws.Cells.Item(822, 1).NumberFormat = "@"
$ws.Cells.Item(822, 1).Value = "2026/02/15"
$ws.Cells.Item(822, 1).Style = "Normal"
$ws.Cells.Item(822, 2).Value = "日"
$ws.Cells.Item(822, 3).Value = 4
$ws.Cells.Item(822, 4).Value = 201
